$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "Hawailian"
$ws.Range("B19").Value = 30
$ws.Range("C19").Value = 7624.5

$ws.Range("A20").Value = "Iced Cappucino"
$ws.Range("B20").Value = 15
$ws.Range("C20").Value = 828

$ws.Range("A21").Value = "Iced Amercano"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 49
